# Capitalize the hex digits (after "0x") in the doip/uds byte-code strings
# found in columns G and H, for easier lookup.
# e.g. "0x02:0xfd:0x00" -> "0x02:0xFD:0x00"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    foreach ($col in @(7, 8)) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2

        if ($val -ne $null -and $val -is [string] -and $val.StartsWith("0x")) {
            $parts = $val -split ":"
            $newParts = @()
            foreach ($p in $parts) {
                if ($p -match "^0x[0-9a-fA-F]{2}$") {
                    $newParts += $p.Substring(0, 2) + $p.Substring(2).ToUpper()
                } else {
                    $newParts += $p
                }
            }
            $newVal = $newParts -join ":"

            if (-not $newVal.Equals($val)) {
                $cell.Value2 = $newVal
            }
        }
    }
}
